$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Octubre de 2020 a las 05:18"

# Row 26 (Pais id 30)
$ws.Range("B26").Value = 319317
$ws.Range("C26").Value = 385
$ws.Range("D26").Value = 304185
$ws.Range("E26").Value = 8552
$ws.Range("G26").Value = 10
$ws.Range("H26").Value = 6580

# Row 31 (Pais id 35)
$ws.Range("B31").Value = 162258
$ws.Range("C31").Value = 5327
$ws.Range("D31").Value = 20272
$ws.Range("E31").Value = 131795
$ws.Range("G31").Value = 16
$ws.Range("H31").Value = 10191

# Row 42 (Pais id 46)
$ws.Range("B42").Value = 108831
$ws.Range("C42").Value = 74
$ws.Range("D42").Value = 104041
$ws.Range("E42").Value = 3044

# Row 153 (Pais id 157)
$ws.Range("B153").Value = 2531
$ws.Range("C153").Value = 35
$ws.Range("D153").Value = 1548
$ws.Range("E153").Value = 946
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 37

# Row 172 (Pais id 176)
$ws.Range("B172").Value = 710
$ws.Range("C172").Value = 3
$ws.Range("E172").Value = 66
